# This script zeroes out a specific set of previously non-zero probability
# cells in the single worksheet of the workbook, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cellsToZero = @(
    "J2", "P2", "S2",
    "J3", "P3",
    "J4", "P4",
    "F6", "J6", "O6", "Q6", "R6", "S6",
    "D10", "F10", "J10", "O10", "Q10", "R10", "S10",
    "J11",
    "J15", "K15", "O15", "S15",
    "F16", "J16", "K16", "S16",
    "F17", "J17", "K17", "O17", "S17",
    "J18", "K18", "O18", "S18",
    "F19", "J19", "K19", "O19", "S19"
)

foreach ($cellRef in $cellsToZero) {
    $ws.Range($cellRef).Value = 0
}
